# Case_2_189 (380 kV case): updated loading_percent results.
# Columns D, J, M, N remain 0 and are left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("B", "C", "E", "F", "G", "H", "I", "K", "L", "O")

$values = @(11.04316897735902, 9.529429008432901, 12.16791198099106, 16.86991607391245, 27.13523703743409, 14.09162271049592, 22.10366211504292, 8.213575834162331, 9.586326816507942, 21.19383117810419)
for ($i = 0; $i -lt $columns.Count; $i++) { $ws.Range($columns[$i] + "2").Value = $values[$i] }

$values = @(10.6995885704849, 9.540705181381357, 12.18492036722794, 15.89584955866815, 27.32840882419302, 14.15005147322406, 22.2231134193017, 7.920408900012639, 9.56803872377392, 21.30633428474998)
for ($i = 0; $i -lt $columns.Count; $i++) { $ws.Range($columns[$i] + "3").Value = $values[$i] }

$values = @(10.48361665991746, 9.54847954306752, 12.19776890279291, 15.26997757108489, 27.45546099089296, 14.18799800524847, 22.30057657730751, 7.752130752249724, 9.558358605380223, 21.37965994831082)
for ($i = 0; $i -lt $columns.Count; $i++) { $ws.Range($columns[$i] + "4").Value = $values[$i] }

$values = @(10.39447099811017, 9.551862181734789, 12.20360963937804, 15.008197319934, 27.50935391524095, 14.2039831692937, 22.33318095241362, 7.686700183262405, 9.554806492871132, 21.41060959123456)
for ($i = 0; $i -lt $columns.Count; $i++) { $ws.Range($columns[$i] + "5").Value = $values[$i] }

$values = @(10.37960370474882, 9.55243684031656, 12.20461601829741, 14.96433081551589, 27.51843057665337, 14.20666902530183, 22.33865759496475, 7.675783322455717, 9.55424046660373, 21.41581332879068)
for ($i = 0; $i -lt $columns.Count; $i++) { $ws.Range($columns[$i] + "6").Value = $values[$i] }

$values = @(10.48241883486052, 9.54852429313137, 12.19784522406672, 15.26647399323133, 27.45617924016707, 14.18821147348513, 22.30101208831493, 7.751251889099764, 9.558309106719548, 21.38007301753582)
for ($i = 0; $i -lt $columns.Count; $i++) { $ws.Range($columns[$i] + "7").Value = $values[$i] }

$values = @(10.92581461135019, 9.533140819660826, 12.17327730461289, 16.5399640634477, 27.20008700334839, 14.11133962519122, 22.14399487402414, 8.113968029699608, 9.579701587654904, 21.2317408521174)
for ($i = 0; $i -lt $columns.Count; $i++) { $ws.Range($columns[$i] + "8").Value = $values[$i] }

$values = @(11.75056248895487, 9.509697990146842, 12.1441785454875, 19.00274580682531, 26.76512414716085, 13.97698730285663, 21.86869778905387, 8.804397219706457, 9.633792084564174, 20.97454426435046)
for ($i = 0; $i -lt $columns.Count; $i++) { $ws.Range($columns[$i] + "9").Value = $values[$i] }

$values = @(12.32326216467891, 9.496538584056063, 12.13441390799951, 20.67494806633232, 26.486872002362, 13.88821680958858, 21.68621524692053, 9.273093751893942, 9.680716621107166, 20.80607928681399)
for ($i = 0; $i -lt $columns.Count; $i++) { $ws.Range($columns[$i] + "10").Value = $values[$i] }

$values = @(12.57551472331739, 9.491426893851381, 12.13248648583461, 21.3917225636224, 26.36933312063486, 13.84997922394826, 21.60747247785345, 9.477370741204718, 9.703572084106501, 20.73388433580428)
for ($i = 0; $i -lt $columns.Count; $i++) { $ws.Range($columns[$i] + "11").Value = $values[$i] }

$values = @(12.66977033574105, 9.489616346580739, 12.1321172690557, 21.65686569030329, 26.32613068858781, 13.83580721508231, 21.57826712289267, 9.553400678062498, 9.712439071444511, 20.70718427041645)
for ($i = 0; $i -lt $columns.Count; $i++) { $ws.Range($columns[$i] + "12").Value = $values[$i] }

$values = @(12.64952819527904, 9.490000724814864, 12.13218076317439, 21.60004134736742, 26.3353768609088, 13.83884573430057, 21.58452978046857, 9.537085800326967, 9.71052005381449, 20.71290620420204)
for ($i = 0; $i -lt $columns.Count; $i++) { $ws.Range($columns[$i] + "13").Value = $values[$i] }

$values = @(12.58329499471333, 9.491275434687401, 12.13244888810653, 21.4136618050453, 26.36575259495346, 13.84880712043312, 21.60505746007557, 9.483652543341419, 9.704297354305151, 20.7316749064925)
for ($i = 0; $i -lt $columns.Count; $i++) { $ws.Range($columns[$i] + "14").Value = $values[$i] }

$values = @(12.54255809686299, 9.492072509327059, 12.13266005938777, 21.29868154950795, 26.38452903298713, 13.85494881430199, 21.617711031808, 9.450749428693747, 9.700513245438525, 20.74325444486031)
for ($i = 0; $i -lt $columns.Count; $i++) { $ws.Range($columns[$i] + "15").Value = $values[$i] }

$values = @(12.30660359926344, 9.496890190706111, 12.13459039351085, 20.62722412089977, 26.49473594643109, 13.89075882370862, 21.69144708754297, 9.259560383480467, 9.679252918986073, 20.81088675523885)
for ($i = 0; $i -lt $columns.Count; $i++) { $ws.Range($columns[$i] + "16").Value = $values[$i] }

$values = @(12.15967672002495, 9.500069243055592, 12.13641802697678, 20.20408069597325, 26.56466454574078, 13.91327590845592, 21.73777446653091, 9.139952964920163, 9.666593541037287, 20.85351434427373)
for ($i = 0; $i -lt $columns.Count; $i++) { $ws.Range($columns[$i] + "17").Value = $values[$i] }

$values = @(12.07439378801057, 9.501980108424439, 12.13770598394411, 19.95656407809801, 26.60573602047484, 13.92642903348634, 21.76482267466392, 9.070318937618136, 9.659454570113427, 20.87845057118777)
for ($i = 0; $i -lt $columns.Count; $i++) { $ws.Range($columns[$i] + "18").Value = $values[$i] }

$values = @(12.04538794025448, 9.502641258686005, 12.13818275146214, 19.87204792380568, 26.6197879578039, 13.93091715114164, 21.77404979915036, 9.0465992766953, 9.657062034603287, 20.88696533483152)
for ($i = 0; $i -lt $columns.Count; $i++) { $ws.Range($columns[$i] + "19").Value = $values[$i] }

$values = @(12.17539808164341, 9.499722308214684, 12.13619897564016, 20.24955283636154, 26.55713246329962, 13.91085803589973, 21.73280124848449, 9.152772523005556, 9.667926453602679, 20.84893330671792)
for ($i = 0; $i -lt $columns.Count; $i++) { $ws.Range($columns[$i] + "20").Value = $values[$i] }

$values = @(12.60278425161968, 9.490897630300346, 12.13236035362476, 21.46857628470577, 26.356794976376, 13.84587287309901, 21.59901135981487, 9.499383449973662, 9.706119395488894, 20.72614475173409)
for ($i = 0; $i -lt $columns.Count; $i++) { $ws.Range($columns[$i] + "21").Value = $values[$i] }

$values = @(12.87469097447437, 9.485859301189711, 12.13195320081224, 22.22866616901552, 26.23348443839835, 13.80519474011153, 21.51514345005809, 9.71817373584385, 9.732314635623601, 20.64961762021699)
for ($i = 0; $i -lt $columns.Count; $i++) { $ws.Range($columns[$i] + "22").Value = $values[$i] }

$values = @(12.73027113044791, 9.488481848064842, 12.1319785678011, 21.82633154458858, 26.2985978110718, 13.82674153597428, 21.55957889459373, 9.602121409205802, 9.718222522538198, 20.69012097149927)
for ($i = 0; $i -lt $columns.Count; $i++) { $ws.Range($columns[$i] + "23").Value = $values[$i] }

$values = @(12.16829297949656, 9.499878898298702, 12.13629726975641, 20.22900810905287, 26.56053501271761, 13.91195050889602, 21.73504835075801, 9.146979507949759, 9.667323410266931, 20.85100305731368)
for ($i = 0; $i -lt $columns.Count; $i++) { $ws.Range($columns[$i] + "24").Value = $values[$i] }

$values = @(11.53289115673614, 9.515323633790306, 12.15000902618814, 18.34778573295695, 26.87556217778073, 14.01158398043651, 21.93969236925276, 8.624173744525544, 9.617880765444999, 21.04052080018431)
for ($i = 0; $i -lt $columns.Count; $i++) { $ws.Range($columns[$i] + "25").Value = $values[$i] }
